$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# Update bladeLengthWheelAnchor value in B17
$ws.Range("B17").Value = 20

# Move the active selection to B17 (matches the author's last cursor position)
$ws.Range("B17").Select()
